$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 33 with the new Mac-Address / Document Type test data record
$ws.Range("A33").Value = 10002
$ws.Range("B33").Value = 110032
$ws.Range("C33").Value = 10032
$ws.Range("D33").Value = "eng"
$ws.Range("E33").Value = $true
$ws.Range("F33").Value = "superadmin"
$ws.Range("G33").Value = "now()"
$ws.Range("H33").Value = "now()"

# Selection / view state change (topLeftCell reset to default, selection moves to E31)
[void]$ws.Range("E31").Select()

# Set calculation mode to manual (calcPr calcMode="manual")
$excel.Calculation = -4135
